# Added exception catch for when output excel file is open
try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # The two tasks on the former "critical" edges are no longer flagged
    # once the summary row below spells out the critical path explicitly.
    $ws.Range("F3").Value = $false
    $ws.Range("F4").Value = $false

    # Append a new summary row describing the critical path.
    # Copy A4's formatting (bold / bordered / centered header style) onto
    # A5 before writing its value so the new row matches the existing
    # "index" column look.
    $ws.Range("A4").Copy()
    $ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A5").Value = 3

    $ws.Range("B5").Value = "Critical Path:"
    $ws.Range("C5").Value = "B -> C"

    # Touch the trailing cells so the row is fully materialized through
    # column F (kept blank, matching D/E/F being empty on this row).
    $ws.Range("D5").Style = "Normal"
    $ws.Range("E5").Style = "Normal"
    $ws.Range("F5").Style = "Normal"
}
catch {
    Write-Host "Could not update output Excel file (it may be open): $_"
}
